# The sheet stored the "flag" columns (B:F) for rows 3-9 as OOXML boolean
# cells (t="b", value 0/1 => FALSE/TRUE). The new format stores the same
# data as plain numbers (no boolean type, value 0/1) instead.
#
# Re-assigning a numeric 0 to each of these cells makes Excel drop the
# boolean type and store them as ordinary numeric cells, matching the
# target "new format for excel".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3:F9").Value = 0
